# The document has three inline pictures (Pearson logo x2 in the footers,
# BTEC logo x1 in the "first page" header) whose non-visual drawing name
# ("wp:docPr/@name", mirrored onto "pic:cNvPr/@name") needs to be swapped:
#   footer (first page)   : image2.png -> image1.png
#   footer (default/odd)  : image2.png -> image1.png
#   header (first page)   : image1.jpg -> image2.jpg
#
# InlineShape.Name is not directly settable in this host, but going via
# Range.ShapeRange(1).Name does reach the drawing's name attribute, so we
# use that indirection for every inline picture that needs renaming.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShape($inlineShape, $newName) {
    $shapeRange = $inlineShape.Range.ShapeRange
    $shapeRange.Item(1).Name = $newName
}

# --- Footers: both the "first page" footer and the default/odd footer ---
# carry the Pearson Edexcel logo (image2.png -> image1.png).
for ($i = 1; $i -le 2; $i++) {
    $footer = $sec.Footers.Item($i)
    if ($footer.Exists) {
        $inlineShapes = $footer.Range.InlineShapes
        for ($j = 1; $j -le $inlineShapes.Count; $j++) {
            Rename-InlineShape $inlineShapes.Item($j) "image1.png"
        }
    }
}

# --- Headers: the "first page" header carries the BTEC logo
# (image1.jpg -> image2.jpg). The default header has no picture.
for ($i = 1; $i -le 2; $i++) {
    $header = $sec.Headers.Item($i)
    if ($header.Exists) {
        $inlineShapes = $header.Range.InlineShapes
        for ($j = 1; $j -le $inlineShapes.Count; $j++) {
            Rename-InlineShape $inlineShapes.Item($j) "image2.jpg"
        }
    }
}
